# CHE_trd_diesel.xlsx edit script
# Sector PASSENGER completed - Added soft-mobility - Improved vehicle speed -
# Long and short travel - Added "choke" constraints
#
# Functionally: replace the sparse World-Bank-sourced cost_import price
# series (rows 76-105, every other year) with a complete SFOE
# (Swiss Federal Office of Energy) sourced series covering every year
# 1990-2019, in CHF2010/litre_diesel, referencing "Table 39" of the
# Gesamtenergiestatistik.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: no content change needed (string indices only shift because
#     unused shared strings are pruned on save) ---

# --- New SFOE-sourced data rows 76-105 -------------------------------------------------

$values = @{
    76 = 0.88
    77 = 0.88
    78 = 0.83
    79 = 0.98
    80 = 0.98
    81 = 0.95
    82 = 1.04
    83 = 1.07
    84 = 1.01
    85 = 1.07
    86 = 1.25
    87 = 1.2
    88 = 1.1499999999999999
    89 = 1.19
    90 = 1.28
    91 = 1.45
    92 = 1.5
    93 = 1.48
    94 = 1.66
    95 = 1.33
    96 = 1.42
    97 = 1.54
    98 = 1.62
    99 = 1.57
    100 = 1.52
    101 = 1.34
    102 = 1.25
    103 = 1.33
    104 = 1.46
    105 = 1.42
}

$unit = "CHF2010/litre_diesel"
$reference = "SFOE"
$link = "https://www.bfe.admin.ch/bfe/de/home/versorgung/statistik-und-geodaten/energiestatistiken/gesamtenergiestatistik.html/"
$note = "Table 39"

foreach ($row in 76..105) {
    $ws.Cells.Item($row, 7).Value = $values[$row]          # G: Value
    $ws.Cells.Item($row, 8).Value = $unit                  # H: Unit
    $ws.Cells.Item($row, 10).Value = $reference             # J: Reference
    $ws.Cells.Item($row, 11).Value = $link                  # K: Link
    $ws.Cells.Item($row, 11).Style = "Hyperlink"
    $ws.Cells.Item($row, 12).Value = $note                  # L: Note
}

# --- Row 106: stray styled (Hyperlink-style) empty cell under K105 ---------------------
$ws.Cells.Item(106, 11).Style = "Hyperlink"

# --- Hyperlinks clean-up: remove the 11 old World-Bank hyperlinks, add one on K74 -------
foreach ($row in @(81, 84, 86, 88, 90, 92, 94, 96, 98, 100, 102)) {
    $ws.Cells.Item($row, 11).Hyperlinks.Delete()
}

$wb.Worksheets.Item(1).Hyperlinks.Add($ws.Cells.Item(74, 11), "https://data.worldbank.org/indicator/EP.PMP.SGAS.CD") | Out-Null

# --- Sheet view: scrolled down to the new data, selection moved -------------------------
$ws.Application.ActiveWindow.ScrollRow = 69
$ws.Range("G106").Select() | Out-Null
